# Insert a new weekly price record for Mango (Vega Central Mapocho de Santiago)
# at row 163, pushing the existing rows 163-252 down to 164-253.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 163..252 down to 164..253, leaving row 163 free for the new record.
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new observation.
$ws.Range("A163").Value = 9
$ws.Range("B163").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C163").Value = "Metropolitana"
$ws.Range("D163").Value = 44455
$ws.Range("E163").Value = 13
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100108
$ws.Range("H163").Value = "Tropicales y subtropicales"
$ws.Range("I163").Value = 100108002
$ws.Range("J163").Value = "Mango"
$ws.Range("K163").Value = "Sin especificar"
$ws.Range("L163").Value = "Primera"
$ws.Range("M163").Value = 630
$ws.Range("N163").Value = 7500
$ws.Range("O163").Value = 8000
$ws.Range("P163").Value = 7722
$ws.Range("Q163").Value = "`$/bandeja 4 kilos"
$ws.Range("R163").Value = "Brasil"
$ws.Range("S163").Value = 1930
$ws.Range("T163").Value = 4
